$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the stray column D (the D6 cell with a leftover fill/border style) ---
$ws.Range("D:D").Delete() | Out-Null

# --- Strip the body-row highlight style (s="2") from the existing data rows 2-6, ---
# --- leaving them with the default/no explicit style, matching the target file. ---
for ($r = 2; $r -le 6; $r++) {
    $ws.Rows.Item($r).ClearFormats()
}

# --- Append the new login-view translation rows (Keys / EN / ZH) ---
$ws.Range("A7").Value2 = 'login_user_placeholder'
$ws.Range("B7").Value2 = 'Please input the username'
$ws.Range("C7").Value2 = '请输入用户名'
$ws.Range("A8").Value2 = 'login_password_placeholder'
$ws.Range("B8").Value2 = 'Please input the password'
$ws.Range("C8").Value2 = '请输入密码'
$ws.Range("A9").Value2 = 'login_verify_code_placeholder'
$ws.Range("B9").Value2 = 'Please input the verify code'
$ws.Range("C9").Value2 = '请输入验证码'
$ws.Range("A10").Value2 = 'login_msg_user_req'
$ws.Range("B10").Value2 = 'Please input the username!'
$ws.Range("C10").Value2 = '请输入用户名！'
$ws.Range("A11").Value2 = 'login_msg_password_req'
$ws.Range("B11").Value2 = 'Please input the password!'
$ws.Range("C11").Value2 = '请输入密码！'
$ws.Range("A12").Value2 = 'login_msg_code_req'
$ws.Range("B12").Value2 = 'Please input the verify code!'
$ws.Range("C12").Value2 = '请输入验证码！'
$ws.Range("A13").Value2 = 'login_msg_mobile_req'
$ws.Range("B13").Value2 = 'Please input mobile number'
$ws.Range("C13").Value2 = '请输入手机号'
$ws.Range("A14").Value2 = 'login_msg_code_sending'
$ws.Range("B14").Value2 = 'Verify code is sending…'
$ws.Range("C14").Value2 = '验证码发送中…'
$ws.Range("A15").Value2 = 'login_msg_code_error'
$ws.Range("B15").Value2 = 'Wrong verify code!'
$ws.Range("C15").Value2 = '您输入的验证码不正确！'
$ws.Range("A16").Value2 = 'login_msg_mobile_error'
$ws.Range("B16").Value2 = 'Wrong mobile number!'
$ws.Range("C16").Value2 = '您的手机号码格式不正确！'
$ws.Range("A17").Value2 = 'login_welcome_title'
$ws.Range("B17").Value2 = 'Welcome'
$ws.Range("C17").Value2 = '欢迎'
$ws.Range("A18").Value2 = 'login_welcome_content'
$ws.Range("B18").Value2 = '{time}, welcome back'
$ws.Range("C18").Value2 = '{time}，欢迎回来'
$ws.Range("A19").Value2 = 'login_msg_login_error'
$ws.Range("B19").Value2 = 'Failed to login'
$ws.Range("C19").Value2 = '登录失败'

# --- Match the saved selection / active cell from the authored workbook ---
$ws.Range("D19").Select() | Out-Null
